$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.117.52"
$ws.Range("E2").Value = "  -1.84%  "
$ws.Range("D3").Value = "1.795.85"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'222.61"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "'0.550"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'32.11"
$ws.Range("E8").Value = "  -2.18%  "
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("D10").Value = "'0.0716"
$ws.Range("E10").Value = "  +3.83%  "
$ws.Range("D11").Value = "'0.0922"
$ws.Range("D12").Value = "2.053.70"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").Value = "1.806.46"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("E14").Value = "  -3.10%  "
$ws.Range("D15").Value = "'0.629"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "34.089.48"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("D18").Value = "'67.99"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").Value = "'245.84"
$ws.Range("E19").Value = "  -3.07%  "
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("E22").Value = "  +2.62%  "
$ws.Range("E23").Value = "  -2.57%  "
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").Value = "'158.88"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "'16.45"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").Value = "'7.03"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("E28").Value = "  -1.93%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'0.0517"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'3.70"
$ws.Range("E31").Value = "  -1.30%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.21"
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("D33").Value = "'3.50"
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("D35").Value = "1.413.85"
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("D36").Value = "'0.643"
$ws.Range("E36").Value = "  +1.72%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("D39").Value = "'0.943"
$ws.Range("E39").Value = "  +4.14%  "
$ws.Range("D40").Value = "'80.07"
$ws.Range("E40").Value = "  -3.52%  "
$ws.Range("E41").Value = "  -2.62%  "
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E43").Value = "  +2.96%  "
$ws.Range("D44").Value = "'5.95"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "'0.0495"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("D46").Value = "1.952.65"
$ws.Range("E47").Value = "  -2.67%  "
$ws.Range("D48").Value = "'105.97"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "'11.87"
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("E51").Value = "  +0.00%  "
